# Añadir funcionalidad de consulta de tareas: se agregan registros de ejemplo
# en la hoja "tareas" para poder probar la consulta (CRUD - Read).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tareas")

# Fila 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "sdad"
$ws.Cells.Item(2, 3).Value = "sadsa"
$ws.Cells.Item(2, 4).Value = "pendiente"
$ws.Cells.Item(2, 5).Value = "dsadas"
$ws.Cells.Item(2, 6).Value = "sadsa"

# Fila 3
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "sdas"
$ws.Cells.Item(3, 3).Value = "sadsa"
$ws.Cells.Item(3, 4).Value = "eje"
$ws.Cells.Item(3, 5).Value = "sdsa"
$ws.Cells.Item(3, 6).Value = "sdsa"

[void]$ws.Range("E5").Select()
